# Movie Randomizer list: merge per-item bulleted paragraphs into one
# comma-separated, quoted list in a single plain paragraph.

$d = $word.ActiveDocument

$titles = @(
    "Matrix",
    "Spider-Man",
    "The Amazing Spider-Man",
    "Interstellar",
    "The Martian",
    "The Disaster Artist",
    "Despicable Me",
    "Casablanca",
    "Argo",
    "Crash",
    "Capote",
    "Mission Impossible: Rogue Nation",
    "Mission Impossible 2",
    "Sonic the Hedgehog",
    "Kingsman: The Secret Service",
    "Iron Man",
    "The Avengers",
    "Star Wars: Empire Strikes Back",
    "Blade Runner",
    "Dances With Wolves",
    "A Clockwork Orange",
    "Sharknado",
    "The Birds",
    "Silence of the Lambs",
    "Misery",
    "Psycho",
    "The Shining",
    "Lord of the Rings: Fellowship of the Ring",
    "Transporter",
    "Die Hard",
    "Galaxy Quest",
    "Toy Story",
    "Ocean’s Eleven",
    "Goldfinger",
    "Get Smart",
    "Casino Royale",
    "Knives Out",
    "Flubber",
    "Hook",
    "Jumanji",
    "Mrs. Doubtfire",
    "El Dorado",
    "Shrek",
    "Tarzan",
    "The Mummy",
    "Space Balls",
    "Robin Hood: Men in Tights",
    "The Evil Dead",
    "Get Out",
    "Nightmare on Elm Street",
    "Halloween",
    "Hocus Pocus",
    "The Greatest Showman",
    "High School Musical",
    "Twilight",
    "Harry Potter and the Sorcerer’s Stone",
    "Harry Potter and the Goblet of Fire",
    "The Hunger Games",
    "Silver Linings Playbook",
    "Full Metal Jacket",
    "The Room",
    "Mortal Kombat",
    "Scary Movie",
    "Speed",
    "Point Break",
    "Doctor Sleep",
    "Batman Begins",
    "The Dark Knight",
    "Batman Returns",
    "She’s All That",
    "Not Another Teen Movie",
    "Godzilla",
    "Rocky"
)

# Locate the first ("Matrix") and last ("Rocky") paragraphs of the list
# by scanning for their distinctive, list-numbered paragraph text.
$startIndex = -1
$endIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($startIndex -eq -1 -and $t -eq "Matrix`r") {
        $startIndex = $i
    }
    if ($startIndex -ne -1 -and $t -eq "Rocky`r") {
        $endIndex = $i
    }
}

$pStart = $d.Paragraphs.Item($startIndex)
$pEnd = $d.Paragraphs.Item($endIndex)

# Merge all the in-between paragraph marks so every title ends up
# back-to-back inside a single paragraph (up to, but not including,
# the very last mark so the final paragraph absorbs "Rocky" too).
$mergeRange = $d.Range($pStart.Range.Start, $d.Paragraphs.Item($endIndex - 1).Range.End)
$mergeRange.Find.Execute("^p", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

$listPara = $d.Paragraphs.Item($startIndex)
$listPara.Range.ListFormat.RemoveNumbers()
$listPara.Style = "Normal"

# Walk the titles in order, quoting each one and stitching
# `", "`-style separators between them.
$cursor = $listPara.Range.Start
$count = $titles.Length
for ($k = 0; $k -lt $count; $k++) {
    $title = $titles[$k]
    $paraEnd = $listPara.Range.End
    $searchRange = $d.Range($cursor, $paraEnd)
    $found = $searchRange.Find.Execute($title, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $title (k=$k)"
    }
    if ($k -eq 0) {
        $searchRange.InsertBefore("“")
    } else {
        $searchRange.InsertBefore("”, “")
    }
    if ($k -eq ($count - 1)) {
        $searchRange.InsertAfter("”")
    }
    $cursor = $searchRange.End
}

Write-Output $listPara.Range.Text
